$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "got dictionary" -> "get a dictionary" (row 6, inference_int8_input_dict)
$ws.Range("B6").Value = 'Same as "inference_int8", but get a dictionary input which has a form of {node name: node input} instead.
This can be used for models with multiple inputs.'

# Add line break (space -> newline) in description texts, and set wrap/row height
$ws.Range("B4").Value = "Floating inference with the Mobilint IR.
This can be used to check the built IR returns the same output as the model."

$ws.Range("B7").Value = "Return the number of add/multiplication operations in the build Mobilint IR.
This can be reduced in later optimization steps."

$ws.Range("B5").Value = "Integer inference with the compiled and quantized model.
The model must be compiled before executing this function."

# Apply wrap text styling to B4:B7 (matches style index 1 used in target)
$ws.Range("B4:B7").WrapText = $true

# Row heights to match new wrapped content
$ws.Rows.Item(4).RowHeight = 33
$ws.Rows.Item(5).RowHeight = 33
$ws.Rows.Item(6).RowHeight = 49.5
$ws.Rows.Item(7).RowHeight = 33

# Update selection to B8
$ws.Range("B8").Select()
